$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204, shifting existing rows 204-236 down to 205-237
$ws.Rows(204).Insert()

# Populate the new row 204 with a new price record (same shape as the surrounding rows)
$ws.Cells.Item(204, 1).Value = 1
$ws.Cells.Item(204, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(204, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(204, 4).Value = 44694
$ws.Cells.Item(204, 5).Value = 15
$ws.Cells.Item(204, 6).Value = "Fruta"
$ws.Cells.Item(204, 7).Value = 100108
$ws.Cells.Item(204, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(204, 9).Value = 100108006
$ws.Cells.Item(204, 10).Value = "Plátano"
$ws.Cells.Item(204, 11).Value = "Sin especificar"
$ws.Cells.Item(204, 12).Value = "Pintón"
$ws.Cells.Item(204, 13).Value = 120
$ws.Cells.Item(204, 14).Value = 16000
$ws.Cells.Item(204, 15).Value = 17000
$ws.Cells.Item(204, 16).Value = 16500
$ws.Cells.Item(204, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(204, 18).Value = "Ecuador"
$ws.Cells.Item(204, 19).Value = 825
$ws.Cells.Item(204, 20).Value = 20
